$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-06-09 Monday" "2025-06-10 Tuesday"

Replace-Text "813÷9=" "496÷5="
Replace-Text "258÷6=" "244÷4="
Replace-Text "406÷8=" "479÷9="
Replace-Text "540÷5=" "733÷4="
Replace-Text "973÷9=" "396÷8="
Replace-Text "389÷7=" "345÷8="
Replace-Text "509÷6=" "694÷3="
Replace-Text "735÷2=" "375÷2="
Replace-Text "440÷7=" "593÷8="
Replace-Text "204÷8=" "666÷9="
Replace-Text "101÷8=" "251÷4="
Replace-Text "858÷3=" "148÷3="
Replace-Text "839÷8=" "123÷9="
Replace-Text "819÷2=" "182÷4="
Replace-Text "185÷7=" "466÷9="
Replace-Text "670÷4=" "123÷7="
Replace-Text "282÷5=" "545÷6="
Replace-Text "637÷5=" "811÷4="
Replace-Text "694÷7=" "421÷7="
Replace-Text "784÷9=" "157÷7="
Replace-Text "924÷2=" "329÷7="
Replace-Text "159÷7=" "278÷7="
Replace-Text "206÷2=" "409÷8="
Replace-Text "730÷4=" "231÷7="
Replace-Text "163÷2=" "933÷6="

Write-Output "Done applying replacements"
